$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.095.28'
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = '1.793.79'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '316.98'
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").Value = '0.5373'
$ws.Range("E7").Value = '  -2.59%  '
$ws.Range("D8").Value = '0.3773'
$ws.Range("E8").Value = '  -2.27%  '
$ws.Range("D9").Value = '0.07462'
$ws.Range("E9").Value = '  -1.91%  '
$ws.Range("D10").Value = '41.67'
$ws.Range("E10").Value = '  -1.97%  '
$ws.Range("D11").Value = '1.093'
$ws.Range("E11").Value = '  -3.16%  '
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").Value = '20.58'
$ws.Range("E13").Value = '  -3.07%  '
$ws.Range("D14").Value = '6.101'
$ws.Range("E14").Value = '  -1.43%  '
$ws.Range("D15").Value = '1.784.85'
$ws.Range("E15").Value = '  -1.17%  '
$ws.Range("D16").Value = '7.206'
$ws.Range("E16").Value = '  -3.09%  '
$ws.Range("D17").Value = '89.11'
$ws.Range("E17").Value = '  -3.19%  '
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("D19").Value = '0.06468'
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("E20").Value = '  +0.08%  '
$ws.Range("D21").Value = '17.29'
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("D22").Value = '5.898'
$ws.Range("E22").Value = '  -1.31%  '
$ws.Range("D23").Value = '28.125.29'
$ws.Range("E24").Value = '  -2.31%  '
$ws.Range("D25").Value = '2.096'
$ws.Range("E25").Value = '  -2.16%  '
$ws.Range("D26").Value = '154.78'
$ws.Range("E26").Value = '  -2.68%  '
$ws.Range("E27").Value = '  -2.37%  '
$ws.Range("D28").Value = '1.990.75'
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("D29").Value = '2.280'
$ws.Range("E29").Value = '  -5.61%  '
$ws.Range("D30").Value = '120.62'
$ws.Range("E30").Value = '  -2.66%  '
$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("D32").Value = '0.1055'
$ws.Range("E32").Value = '  +2.98%  '
$ws.Range("D33").Value = '3.655'
$ws.Range("E33").Value = '  -0.88%  '
$ws.Range("D34").Value = '5.555'
$ws.Range("E34").Value = '  -3.80%  '
$ws.Range("D35").Value = '0.06545'
$ws.Range("E35").Value = '  +1.73%  '
$ws.Range("D36").Value = '0.2259'
$ws.Range("E36").Value = '  -2.47%  '
$ws.Range("E37").Value = '  -1.94%  '
$ws.Range("D38").Value = '5.017'
$ws.Range("E38").Value = '  -3.15%  '
$ws.Range("D39").Value = '8.443'
$ws.Range("E39").Value = '  -4.17%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.6168'
$ws.Range("E40").Value = '  -3.94%  '
$ws.Range("B41").Value = 'WEMIXTOKEN'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").Value = '1.446'
$ws.Range("E41").Value = '  +4.49%  '
$ws.Range("D42").Value = '11.05'
$ws.Range("E42").Value = '  -5.14%  '
$ws.Range("E43").Value = '  +0.79%  '
$ws.Range("D44").Value = '0.9999'
$ws.Range("D45").Value = '13.31'
$ws.Range("E45").Value = '  -2.40%  '
$ws.Range("D46").Value = '3.673'
$ws.Range("E46").Value = '  -0.26%  '
$ws.Range("D47").Value = '0.5777'
$ws.Range("E47").Value = '  -3.55%  '
$ws.Range("D48").Value = '126.32'
$ws.Range("E48").Value = '  -0.65%  '
$ws.Range("D49").Value = '1.187'
$ws.Range("E49").Value = '  +3.11%  '
$ws.Range("D50").Value = '1.924'
$ws.Range("E50").Value = '  -3.11%  '
$ws.Range("E51").Value = '  -1.24%  '
